# Sankey data workbook revision - "revised data based on research feedback"
#
# Adds a new column K with ratio/percentage calculations alongside the
# existing Category/Value (I/J) table, and renames the five
# "<Race> Teaching degree" category labels to "<Race>-Teaching degree".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the "Teaching degree" category labels (space -> hyphen)
# ---------------------------------------------------------------------
$ws.Range("I17").Value = "White-Teaching degree"
$ws.Range("I18").Value = "Black-Teaching degree"
$ws.Range("I19").Value = "Hispanic-Teaching degree"
$ws.Range("I20").Value = "Asian-Teaching degree"
$ws.Range("I21").Value = "Other-Teaching degree"

# ---------------------------------------------------------------------
# 2. New column K: "All <Race>" rows (2-6) get a flat 100% marker
# ---------------------------------------------------------------------
foreach ($r in 2..6) {
    $ws.Range("K$r").Value = 1
    $ws.Range("K$r").NumberFormat = "0%"
}

# ---------------------------------------------------------------------
# 3. New column K: "<Race> HS diploma" rows (7-11) -- share of row 2-6
# ---------------------------------------------------------------------
$ws.Range("K7").Formula = "=J7/J2"
$ws.Range("K8").Formula = "=J8/J3"
$ws.Range("K9").Formula = "=J9/J4"
$ws.Range("K10").Formula = "=J10/J5"
$ws.Range("K11").Formula = "=J11/J6"
$ws.Range("K7:K11").Style = "Normal"

# ---------------------------------------------------------------------
# 4. New column K: remaining ratio rows (12-31), formatted "0.00"
# ---------------------------------------------------------------------
$ws.Range("K12").Formula = "=J12/J7"
$ws.Range("K13").Formula = "=J13/J8"
$ws.Range("K14").Formula = "=J14/J9"
$ws.Range("K15").Formula = "=J15/J10"
$ws.Range("K16").Formula = "=J16/J6"
$ws.Range("K17").Formula = "=J17/J12"
$ws.Range("K18").Formula = "=J18/J13"
$ws.Range("K19").Formula = "=J19/J14"
$ws.Range("K20").Formula = "=J20/J15"
$ws.Range("K21").Formula = "=J21/J16"
$ws.Range("K22").Formula = "=J22/J17"
$ws.Range("K23").Formula = "=J23/J18"
$ws.Range("K24").Formula = "=J24/J19"
$ws.Range("K25").Formula = "=J25/J20"
$ws.Range("K26").Formula = "=J26/J21"
$ws.Range("K27").Formula = "=J27/J12"
$ws.Range("K28").Formula = "=J28/J13"
$ws.Range("K29").Formula = "=J29/J14"
$ws.Range("K30").Formula = "=J30/J15"
$ws.Range("K31").Formula = "=J31/J16"
$ws.Range("K12:K31").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 5. New column K: summary rows (32-42) -- same "0.00" style, left blank
# ---------------------------------------------------------------------
foreach ($r in 32..42) {
    $ws.Range("K$r").NumberFormat = "0.00"
}

# ---------------------------------------------------------------------
# 6. Row-height cleanup: rows whose wrapped text now fits on one line
#    (triggered by the column/content edits above) drop their explicit
#    30pt height back to the sheet default.
# ---------------------------------------------------------------------
foreach ($r in @(7, 9, 11, 33, 35, 36, 38, 39, 41, 42)) {
    $ws.Rows("$r").AutoFit()
}

# ---------------------------------------------------------------------
# 7. Leave the selection where the editor's cursor ended up
# ---------------------------------------------------------------------
$ws.Range("K32").Select()
